$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the timestamp column (A2:A27) with the new flexibility-input datetimes
# (Excel serial date-times, kept at full float precision to match source data).
$ws.Range("A2").Value2  = 45431.916666666664
$ws.Range("A3").Value2  = 45431.958333333336
$ws.Range("A4").Value2  = 45432
$ws.Range("A5").Value2  = 45432.041666666664
$ws.Range("A6").Value2  = 45432.08333321759
$ws.Range("A7").Value2  = 45432.124999826388
$ws.Range("A8").Value2  = 45432.166666435187
$ws.Range("A9").Value2  = 45432.208333043978
$ws.Range("A10").Value2 = 45432.249999652777
$ws.Range("A11").Value2 = 45432.291666261575
$ws.Range("A12").Value2 = 45432.333332870374
$ws.Range("A13").Value2 = 45432.374999479165
$ws.Range("A14").Value2 = 45432.416666087964
$ws.Range("A15").Value2 = 45432.458332696762
$ws.Range("A16").Value2 = 45432.499999305554
$ws.Range("A17").Value2 = 45432.541665914352
$ws.Range("A18").Value2 = 45432.583332523151
$ws.Range("A19").Value2 = 45432.624999131942
$ws.Range("A20").Value2 = 45432.66666574074
$ws.Range("A21").Value2 = 45432.708332349539
$ws.Range("A22").Value2 = 45432.74999895833
$ws.Range("A23").Value2 = 45432.791665567129
$ws.Range("A24").Value2 = 45432.833332175927
$ws.Range("A25").Value2 = 45432.874998784719
$ws.Range("A26").Value2 = 45432.916665393517
$ws.Range("A27").Value2 = 45432.958332002316

# Reset the view: scroll back to the top-left (removes the topLeftCell="A15"
# freeze) and select A2:A27 with A2 as the active cell.
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A2:A27").Select() | Out-Null
